$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.190.06'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '1.827.40'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('D4').Value = "'1.003"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').Value = "'233.91"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.18%  '
$ws.Range('D6').Value = "'0.5972"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.42%  '
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').Value = "'0.06948"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.77%  '
$ws.Range('D9').Value = "'0.2742"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.09%  '
$ws.Range('D10').Value = "'23.24"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.35%  '
$ws.Range('D11').Value = "'0.07609"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.35%  '
$ws.Range('D12').Value = '1.831.80'
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('D13').Value = "'4.752"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.19%  '
$ws.Range('D14').Value = "'0.6245"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.04%  '
$ws.Range('D15').Value = "'0.000009674"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.24%  '
$ws.Range('D16').Value = "'78.13"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.40%  '
$ws.Range('D17').Value = '28.844.71'
$ws.Range('E17').Value = '  -1.65%  '
$ws.Range('D18').Value = "'5.701"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -9.12%  '
$ws.Range('D19').Value = "'221.07"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.63%  '
$ws.Range('D20').Value = "'1.005"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('D21').Value = "'11.49"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.56%  '
$ws.Range('D22').Value = "'6.854"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.03%  '
$ws.Range('D23').Value = "'1.005"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.32%  '
$ws.Range('D24').Value = "'155.51"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.06%  '
$ws.Range('D25').Value = "'7.949"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.18%  '
$ws.Range('D26').Value = "'0.1286"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.05%  '
$ws.Range('D27').Value = "'16.50"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.70%  '
$ws.Range('D28').Value = "'0.06637"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -8.44%  '
$ws.Range('E29').Value = '  -3.07%  '
$ws.Range('D30').Value = "'1.439"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.35%  '
$ws.Range('D31').Value = "'3.837"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.64%  '
$ws.Range('D32').Value = "'3.750"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.27%  '
$ws.Range('D33').Value = "'1.089"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.24%  '
$ws.Range('D34').Value = "'1.713"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.64%  '
$ws.Range('D35').Value = "'0.6448"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.93%  '
$ws.Range('D36').Value = "'2.543"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.19%  '
$ws.Range('D37').Value = "'2.730"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.49%  '
$ws.Range('D38').Value = "'0.01731"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.30%  '
$ws.Range('D39').Value = "'6.512"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.79%  '
$ws.Range('D40').Value = '1.173.43'
$ws.Range('E40').Value = '  -4.79%  '
$ws.Range('D41').Value = "'0.8949"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.21%  '
$ws.Range('D42').Value = "'1.005"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.36%  '
$ws.Range('D43').Value = '1.977.06'
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('D44').Value = "'100.46"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('D45').Value = "'61.98"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.96%  '
$ws.Range('D46').Value = "'0.00000000113"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.19%  '
$ws.Range('D47').Value = "'0.05545"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.07%  '
$ws.Range('D48').Value = "'8.415"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.22%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = "'1.576"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.00%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = "'0.4550"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.52%  '
$ws.Range('D51').Value = "'0.3640"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.22%  '
